# Append updated data rows (10/09/2021 - 20/09/2021) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data: row, serial-date, B (nuovi pos.), C (somma mobile 7gg.), D (somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(375, 44449, 1, 1, 48.07692307692308),
    @(376, 44450, 1, 2, 96.15384615384616),
    @(377, 44451, 1, 3, 144.2307692307692),
    @(378, 44452, 3, 6, 288.4615384615385),
    @(379, 44453, 0, 6, 288.4615384615385),
    @(380, 44454, 0, 6, 288.4615384615385),
    @(381, 44455, 0, 6, 288.4615384615385),
    @(382, 44456, 0, 5, 240.3846153846154),
    @(383, 44457, 0, 4, 192.3076923076923),
    @(384, 44458, 0, 3, 144.2307692307692),
    @(385, 44459, 0, 0, 0)
)

foreach ($item in $data) {
    $r = $item[0]
    $dateSerial = $item[1]
    $b = $item[2]
    $c = $item[3]
    $d = $item[4]

    # Column A: date value, formatted/styled like the row above it (r-1)
    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $dateSerial

    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}

$excel.CutCopyMode = $false
